$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update D2: "no" -> "yes"
$ws.Range("D2").Value = "yes"

# Update B4: "DrugsdoDrugsdo" -> "12gupta"
$ws.Range("B4").Value = "12gupta"

# Update D4: "yes" -> "no"
$ws.Range("D4").Value = "no"

# Add new row 5
$ws.Range("A5").Value = "siri@gmail.com"
$ws.Range("B5").Value = "12siri"
$ws.Range("C5").Value = "127.0.0.1"
$ws.Range("D5").Value = "no"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "1"
$ws.Range("F5").Value = "alexa"
